$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): force text storage per cell so values like
#     "247.67" or "0.9990" keep their exact text representation
#     (NumberFormat "@" before the write, ClearFormats after to drop the
#     temporary style back to the sheet default -- applied per-cell since
#     batching into one multi-area Range leaves some members numeric). ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.866.07"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.809.83"
$ws.Range("D3").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.67"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4955"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2783"
$ws.Range("D8").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.808.60"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.82"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07084"
$ws.Range("D12").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.14"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.695"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.893.75"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9992"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007327"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9989"
$ws.Range("D19").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.049.41"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.588"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.877"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.356"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.33"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "129.33"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.39"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.890"
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.414"
$ws.Range("D29").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.149"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08351"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.822"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04970"
$ws.Range("D33").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6733"
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.687"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.325"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.753"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9531"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.132"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01595"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9990"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4100"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.92"
$ws.Range("D44").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1220"
$ws.Range("D46").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.133"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.72"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3639"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.312"
$ws.Range("D51").ClearFormats()

# --- Coin name (B) / Link (C) swap for rows 48-49 ---
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("B49").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"

# --- Volume(1h) column (E): plain text percentages, safe to set directly ---
$ws.Range("E2").Value = "  +8.07%  "
$ws.Range("E3").Value = "  +4.97%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +2.58%  "
$ws.Range("E8").Value = "  +8.10%  "
$ws.Range("E9").Value = "  +3.89%  "
$ws.Range("E10").Value = "  +5.01%  "
$ws.Range("E11").Value = "  +6.01%  "
$ws.Range("E12").Value = "  +3.67%  "
$ws.Range("E13").Value = "  +7.09%  "
$ws.Range("E14").Value = "  +9.34%  "
$ws.Range("E15").Value = "  +5.39%  "
$ws.Range("E16").Value = "  +8.87%  "
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  +8.24%  "
$ws.Range("E21").Value = "  +5.50%  "
$ws.Range("E23").Value = "  +3.75%  "
$ws.Range("E24").Value = "  +6.08%  "
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("E26").Value = "  +21.54%  "
$ws.Range("E27").Value = "  +7.75%  "
$ws.Range("E28").Value = "  +7.06%  "
$ws.Range("E29").Value = "  +3.39%  "
$ws.Range("E30").Value = "  +3.62%  "
$ws.Range("E31").Value = "  +5.78%  "
$ws.Range("E32").Value = "  +4.50%  "
$ws.Range("E34").Value = "  +9.66%  "
$ws.Range("E35").Value = "  +9.14%  "
$ws.Range("E36").Value = "  +3.56%  "
$ws.Range("E37").Value = "  +15.70%  "
$ws.Range("E38").Value = "  +12.71%  "
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("E40").Value = "  +9.56%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  +7.27%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("E45").Value = "  +5.91%  "
$ws.Range("E46").Value = "  +6.05%  "
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("E49").Value = "  +5.74%  "
$ws.Range("E50").Value = "  +8.91%  "
$ws.Range("E51").Value = "  +5.93%  "

